# Rotate values in columns A, I, Q, R among rows 9, 10, 11:
#   new row 9  <- old row 10
#   new row 10 <- old row 11
#   new row 11 <- old row 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Antal" column (I) stores its numbers as text, so force the
# cell format to Text before writing the value to keep it a string
# instead of Excel auto-converting it to a number.
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I11").NumberFormat = "@"

# Row 9 -> new values (previously in row 10)
$ws.Range("A9").Value = 111675587
$ws.Range("I9").Value = "3"
$ws.Range("Q9").Value = 690344.8588249951
$ws.Range("R9").Value = 6661440.743740954

# Row 10 -> new values (previously in row 11)
$ws.Range("A10").Value = 111675586
$ws.Range("I10").Value = "2"
$ws.Range("Q10").Value = 690348.8581766916
$ws.Range("R10").Value = 6661440.95072202

# Row 11 -> new values (previously in row 9)
$ws.Range("A11").Value = 111675585
$ws.Range("I11").Value = "1"
$ws.Range("Q11").Value = 690349.9096738817
$ws.Range("R11").Value = 6661440.004307052
